# Append Q-Learning / Joueur Humain session rows (9-13) to the results log.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 9;  Date = "2025-06-25 17:52:57"; Agent = "Q Learning";    Env = "{'win': 1, 'loss': 0, 'draw': 1}" },
    @{ Row = 10; Date = "2025-06-25 17:53:04"; Agent = "Q Learning";    Env = "{'win': 1, 'loss': 1, 'draw': 2}" },
    @{ Row = 11; Date = "2025-06-25 17:53:07"; Agent = "Q Learning";    Env = "{'win': 2, 'loss': 2, 'draw': 2}" },
    @{ Row = 12; Date = "2025-06-25 17:53:10"; Agent = "Q Learning";    Env = "{'win': 3, 'loss': 2, 'draw': 3}" },
    @{ Row = 13; Date = "2025-06-25 17:53:44"; Agent = "Joueur Humain"; Env = "{'win': 0, 'loss': 1, 'draw': 1}" }
)

foreach ($r in $rows) {
    $rowIndex = $r.Row
    $ws.Cells.Item($rowIndex, 1).Value = $r.Date
    $ws.Cells.Item($rowIndex, 2).Value = $r.Agent
    $ws.Cells.Item($rowIndex, 3).Value = $r.Env

    # Columns D..I stay empty (inline/empty text cells, same as earlier rows
    # before hyperparameter tracking columns were populated). Forcing a
    # quote-prefixed empty value then stripping the format mirrors Excel's
    # own "typed into an empty text cell" state without leaving a style.
    for ($col = 4; $col -le 9; $col++) {
        $cell = $ws.Cells.Item($rowIndex, $col)
        $cell.Value = "'"
        $cell.ClearFormats()
    }
}

$null = $ws.Range("A1").Select()
